$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2010-18")

# Row 31 is a copy of row 30 (columns A:R), except column B gets a new
# label string ("Baseline 2010-18 C549" instead of "Baseline 2010-18 C502").
for ($i = 1; $i -le 18; $i++) {
    $srcCell = $ws.Cells.Item(30, $i)
    $dstCell = $ws.Cells.Item(31, $i)
    $dstCell.Value = $srcCell.Value2
    if ($i -ge 4) {
        $dstCell.NumberFormat = $srcCell.NumberFormat
    }
}
$ws.Range("B31").Value = "Baseline 2010-18 C549"

# Move the active selection to the cell below the newly added row, matching
# where the cursor lands after entering a new row of data.
$ws.Range("B32").Select() | Out-Null
